$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily format the Price column as Text so numeric-looking
# strings (e.g. "568.41") are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.782.43"
$ws.Range("E2").Value = "  +0.14%  "

$ws.Range("D3").Value = "3.364.46"
$ws.Range("E3").Value = "  -0.59%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "568.41"

$ws.Range("D6").Value = "137.45"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "7.64"
$ws.Range("E9").Value = "  +2.02%  "

$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").Value = "0.380"
$ws.Range("E11").Value = "  -4.64%  "

$ws.Range("D12").Value = "3.941.07"
$ws.Range("E12").Value = "  -0.54%  "

$ws.Range("D13").Value = "0.125"
$ws.Range("E13").Value = "  +1.31%  "

$ws.Range("D14").Value = "27.69"
$ws.Range("E14").Value = "  -2.18%  "

$ws.Range("D15").Value = "3.363.82"
$ws.Range("E15").Value = "  -0.53%  "

$ws.Range("E16").Value = "  -2.05%  "

$ws.Range("D17").Value = "60.930.37"
$ws.Range("E17").Value = "  +0.19%  "

$ws.Range("E18").Value = "  -2.77%  "

$ws.Range("D19").Value = "13.48"
$ws.Range("E19").Value = "  -4.10%  "

$ws.Range("D20").Value = "8.87"
$ws.Range("E20").Value = "  -1.89%  "

$ws.Range("D21").Value = "380.81"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("D22").Value = "75.58"
$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("E23").Value = "  -2.56%  "

$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  -7.03%  "

$ws.Range("E26").Value = "  +6.63%  "

$ws.Range("E27").Value = "  -0.28%  "

$ws.Range("E28").Value = "  -4.38%  "

$ws.Range("E29").Value = "  -2.57%  "

$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("E32").Value = "  -6.51%  "

$ws.Range("D33").Value = "22.89"
$ws.Range("E33").Value = "  -3.51%  "

$ws.Range("D34").Value = "167.34"
$ws.Range("E34").Value = "  +0.82%  "

$ws.Range("D35").Value = "6.80"
$ws.Range("E35").Value = "  -2.74%  "

$ws.Range("E36").Value = "  -2.67%  "

$ws.Range("D37").Value = "3.400.51"
$ws.Range("E37").Value = "  -0.44%  "

$ws.Range("E38").Value = "  -3.67%  "

$ws.Range("D39").Value = "0.0753"
$ws.Range("E39").Value = "  -3.12%  "

$ws.Range("E40").Value = "  -9.94%  "

$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("D42").Value = "4.32"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("E43").Value = "  -3.83%  "

$ws.Range("D44").Value = "2.455.72"
$ws.Range("E44").Value = "  -2.69%  "

$ws.Range("E45").Value = "  -2.76%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("E47").Value = "  -3.57%  "

$ws.Range("D48").Value = "22.06"
$ws.Range("E48").Value = "  -6.86%  "

$ws.Range("E49").Value = "  -5.22%  "

$ws.Range("D50").Value = "1.96"
$ws.Range("E50").Value = "  -5.52%  "

$ws.Range("D51").Value = "0.201"
$ws.Range("E51").Value = "  -3.89%  "

# Restore the original (default) cell style for the Price column
# now that the values are safely stored as text.
$ws.Range("D2:D51").Style = "Normal"
